$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K2: special (non-shared) formula with absolute references
$ws.Range("K2").Formula = '=IF(ISEVEN(ROW($K2)),"writeSPI1("&$C2&");","writeSPI1("&$A1&");")'

# K3:K37: shared-style formula (relative references), row by row
for ($r = 3; $r -le 37; $r++) {
    $prev = $r - 1
    $formula = '=IF(ISEVEN(ROW(K' + $r + ')),"writeSPI1("&C' + $r + '&");","writeSPI1("&A' + $prev + '&");")'
    $ws.Range("K$r").Formula = $formula
}

# O3:O36: shared-style formula (absolute column, relative row)
for ($r = 3; $r -le 36; $r++) {
    $prev = $r - 1
    $formula = '=IF(ISODD(ROW($K' + $r + ')),"writeSPI1("&$C' + $r + '&");","writeSPI1("&$A' + $prev + '&");")'
    $ws.Range("O$r").Formula = $formula
}

# Match the saved selection/active cell state
$ws.Range("K2:K37").Select()
